# Perform Operation Using Pandas
# Insert the DataFrame's default index column ("Unnamed: 0") into column A,
# matching the header style used by the other header cells, and clear the
# old header-style formatting that had been left on the index values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new index-column header.
$ws.Range("A1").Value = "Unnamed: 0"

# Give A1 the same formatting as the rest of the header row (bold, centered,
# bordered) by copying B1's format onto it.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# The index values (A2:A4) no longer carry the header styling.
$ws.Range("A2:A4").ClearFormats()
